$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '60.104.30'
$ws.Range('E2').Value = '  -3.92%  '

# Row 3
$ws.Range('D3').Value = '3.310.65'
$ws.Range('E3').Value = '  -3.87%  '

# Row 4
$ws.Range('E4').Value = '  +0.10%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '559.87'
$ws.Range('E5').Value = '  -3.29%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.43'
$ws.Range('E6').Value = '  -2.35%  '

# Row 7
$ws.Range('E7').Value = '  +0.07%  '

# Row 8
$ws.Range('E8').Value = '  -0.01%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.81'
$ws.Range('E9').Value = '  -2.83%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.120'
$ws.Range('E10').Value = '  -2.53%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.409'
$ws.Range('E11').Value = '  -0.26%  '

# Row 12
$ws.Range('D12').Value = '3.901.37'
$ws.Range('E12').Value = '  -3.27%  '

# Row 13
$ws.Range('E13').Value = '  +0.95%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.38'
$ws.Range('E14').Value = '  -3.21%  '

# Row 15
$ws.Range('D15').Value = '3.347.77'
$ws.Range('E15').Value = '  -2.94%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000166'
$ws.Range('E16').Value = '  -2.90%  '

# Row 17
$ws.Range('D17').Value = '60.226.65'
$ws.Range('E17').Value = '  -3.78%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.19'
$ws.Range('E18').Value = '  -2.44%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.32'
$ws.Range('E19').Value = '  -2.27%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.76'
$ws.Range('E20').Value = '  -2.62%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '373.53'
$ws.Range('E21').Value = '  -3.19%  '

# Row 22
$ws.Range('B22').Value = 'Litecoin'
$ws.Range('C22').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '74.25'

# Row 23
$ws.Range('B23').Value = 'Polygon'
$ws.Range('C23').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.551'
$ws.Range('E23').Value = '  -2.25%  '

# Row 24
$ws.Range('E24').Value = '  -0.01%  '

# Row 25
$ws.Range('D25').Value = '3.486.17'
$ws.Range('E25').Value = '  -2.69%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000106'
$ws.Range('E26').Value = '  -7.45%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.172'
$ws.Range('E27').Value = '  -5.82%  '

# Row 28
$ws.Range('E28').Value = '  +0.32%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.22'
$ws.Range('E29').Value = '  -5.11%  '

# Row 30
$ws.Range('E30').Value = '  +0.12%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.06'
$ws.Range('E31').Value = '  -2.56%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.60'
$ws.Range('E32').Value = '  -5.05%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '22.61'
$ws.Range('E33').Value = '  -2.50%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.28'
$ws.Range('E34').Value = '  -4.08%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.22'
$ws.Range('E35').Value = '  -2.74%  '

# Row 36
$ws.Range('B36').Value = 'Monero'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '166.91'
$ws.Range('E36').Value = '  -1.25%  '

# Row 37
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.53'
$ws.Range('E37').Value = '  -5.91%  '

# Row 38
$ws.Range('B38').Value = 'Aptos'
$ws.Range('C38').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.75'
$ws.Range('E38').Value = '  -2.69%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '27.70'
$ws.Range('E39').Value = '  -13.21%  '

# Row 40
$ws.Range('D40').Value = '3.369.33'
$ws.Range('E40').Value = '  -3.11%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0736'
$ws.Range('E41').Value = '  -4.71%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '41.99'
$ws.Range('E42').Value = '  -1.39%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.752'
$ws.Range('E43').Value = '  -4.16%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.22'
$ws.Range('E44').Value = '  -3.27%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.12'
$ws.Range('E45').Value = '  -4.43%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.59'
$ws.Range('E46').Value = '  -5.64%  '

# Row 47
$ws.Range('D47').Value = '2.431.78'
$ws.Range('E47').Value = '  -5.50%  '

# Row 48
$ws.Range('B48').Value = 'FirstDigitalUSD'
$ws.Range('C48').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.00'
$ws.Range('E48').Value = '  +0.06%  '

# Row 49
$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.61'
$ws.Range('E49').Value = '  -4.08%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '21.84'
$ws.Range('E50').Value = '  -3.31%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0258'
$ws.Range('E51').Value = '  -2.90%  '
